$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the duplicated "Contact" row (was row 11), shifting subsequent rows up by one.
$ws.Rows.Item(11).Delete()

# Update Version value (row 3)
$ws.Range("B3").Value = "6.0.0"

# Update Date value (row 8)
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Update Publisher value (row 9)
$ws.Range("B9").Value = "Alvearie Team"

# Replace the remaining "Contact" row (now row 10) with Jurisdiction info
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"
